$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateTimeFmt = "[$-409]m/d/yy\ h:mm\ AM/PM;@"
$timeFmt     = "h:mm;@"

# ---------------------------------------------------------------------------
# 1. Header row: center every header cell; "Duration" (D1) switches from the
#    generic date/time format to a plain h:mm format.
# ---------------------------------------------------------------------------
$ws.Range("A1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").NumberFormat = $timeFmt
$ws.Range("E1").HorizontalAlignment = -4108
$ws.Range("F1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 2. New data rows under the table.
#    B2 = 8:00 AM (time-of-day serial), C2 = NOW(), D2 = C2-B2 (elapsed time)
# ---------------------------------------------------------------------------
$ws.Range("B2").HorizontalAlignment = -4108
$ws.Range("B2").Value = 0.33333333333333331

$ws.Range("C2").HorizontalAlignment = -4108
$ws.Range("C2").Formula = "=NOW()"

$ws.Range("D2").HorizontalAlignment = -4108
$ws.Range("D2").NumberFormat = $timeFmt
$ws.Range("D2").Formula = "=C2-B2"

# Row that used to hold the placeholder text, now moved down to B3 with real
# (text) date content.
$ws.Range("B3").HorizontalAlignment = -4108
$ws.Range("B3").Value = "6/13/22  8:00 AM"

# B4: a TIME() formula, alongside the existing "Target hrs/day..." note (I4).
$ws.Range("B4").HorizontalAlignment = -4108
$ws.Range("B4").Formula = "=TIME( 1, 2, 3)"

# B5: a plain date/time serial value.
$ws.Range("B5").HorizontalAlignment = -4108
$ws.Range("B5").Value = 36892.041666666664

# ---------------------------------------------------------------------------
# 3. Widen/best-fit the B:C and D columns for the new, wider content.
# ---------------------------------------------------------------------------
$ws.Columns("B:C").ColumnWidth = 19.43
$ws.Columns("D:D").ColumnWidth = 8.71

# ---------------------------------------------------------------------------
# 4. Scratch format probes further down the sheet (format-only cells; no
#    values were ever committed there).
# ---------------------------------------------------------------------------
$ws.Range("E9").NumberFormat = $dateTimeFmt
$ws.Range("F9").NumberFormat = $dateTimeFmt
$ws.Range("G9").NumberFormat = $timeFmt
$ws.Range("H9").NumberFormat = $dateTimeFmt
$ws.Range("I9").NumberFormat = $dateTimeFmt

$ws.Range("E10").NumberFormat = $dateTimeFmt
$ws.Range("F10").NumberFormat = $dateTimeFmt
$ws.Range("G10").NumberFormat = $dateTimeFmt
$ws.Range("H10").NumberFormat = $timeFmt
$ws.Range("I10").NumberFormat = "h:mm AM/PM"

$ws.Range("E11").NumberFormat = $dateTimeFmt
$ws.Range("F11").NumberFormat = $dateTimeFmt
$ws.Range("G11").NumberFormat = $timeFmt
$ws.Range("I11").NumberFormat = "d-mmm-yy"

# ---------------------------------------------------------------------------
# 5. Selection cursor, matching the author's last click before saving.
# ---------------------------------------------------------------------------
$ws.Range("B9").Select() | Out-Null
